# Remove the first list item ("Jennifer Kramer: Surf") entirely,
# including its paragraph mark, leaving only the remaining
# "Swirly Flowers: Backgrounds Etc." list item.
$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq "Jennifer Kramer: Surf") {
        $p.Range.Delete()
        break
    }
}
